# VolSkill.xlsx: adopt the new "standard" layout.
#
# Old layout:
#   A1            = "VolTeer.Vol.tblVolSkill"            (single fixed table name)
#   A2:C2         = "SkillID" | "VolID" | "Query"
#   A3:C7         = SkillID | VolID | INSERT-statement formula (referencing $A$1/$A$2/$B$2)
#
# New layout:
#   B1            = "VolTeer.Vol.tblVolSkill"            (kept, but no longer used by the formula)
#   A2:D2         = "Table" | "SkillID" | "VolID" | "Query"
#   A3:D5         = "Vol.tblVolSkill" (per row) | SkillID | VolID | INSERT-statement formula
#                   (now referencing the row's own column A, and $B$2/$C$2)
#   Rows 6 & 7 (which only had NULL keys) are dropped; the real VolID that used to
#   live in row 6 is folded into row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new column A; everything that used to be in A/B/C shifts to B/C/D.
$ws.Columns.Item(1).Insert()

# 2) Drop the old rows 6 and 7 (both had a NULL key) - delete bottom-up so indices stay valid.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# 3) Populate the new column A: header "Table" and per-row table name "Vol.tblVolSkill".
$ws.Cells.Item(2, 1).Value = "Table"
for ($r = 3; $r -le 5; $r++) {
    $ws.Cells.Item($r, 1).Value = "Vol.tblVolSkill"
}

# 4) Row 5's VolID (column C) used to be NULL; the standard now carries over the
#    real GUID that used to live in the deleted row 6.
$ws.Cells.Item(5, 3).Value = "37CDE19A-E9F2-4AA5-BB7F-64D633D1CE6D"

# 5) Rebuild the INSERT-statement formulas (now in column D) to read the table name
#    from each row's own column A instead of the single fixed header cell.
for ($r = 3; $r -le 5; $r++) {
    $ws.Cells.Item($r, 4).Formula = '=(((((((((("INSERT INTO " & A' + $r + ') &" (") & $B$2) & ",") & $C$2) & ") VALUES(") & B' + $r + ') & ",") & ",") & C' + $r + ') & ")"'
}
